# Weekly update: a new daily price record was inserted ahead of the existing
# row 68, pushing every subsequent record down by one row (68->69, ..., 188->189).
# The new record itself is written into the now-empty row 68.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68; Excel shifts rows 68..188 down to 69..189
# (and the sheet's dimension grows from R188 to R189 automatically).
$ws.Rows.Item(68).Insert()

# Populate the newly-inserted row 68 with the new day's record.
$ws.Cells.Item(68, 1).Value = 4
$ws.Cells.Item(68, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(68, 3).Value = "Los Lagos"
$ws.Cells.Item(68, 4).Value = 44533
$ws.Cells.Item(68, 5).Value = 10
$ws.Cells.Item(68, 6).Value = 100112043
$ws.Cells.Item(68, 7).Value = "Pepino ensalada"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 400
$ws.Cells.Item(68, 11).Value = 11000
$ws.Cells.Item(68, 12).Value = 11000
$ws.Cells.Item(68, 13).Value = 11000
$ws.Cells.Item(68, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 183
$ws.Cells.Item(68, 17).Value = 60
$ws.Cells.Item(68, 18).Value = "Hortaliza"
